# Update cryptos price/volume data (values scraped on 2024-01-01)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.064.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.310.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.02"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.995"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.664.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.309.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.957.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000105"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0868"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("E34").Value = "  +6.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0356"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.231"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.677.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
